$d = $word.ActiveDocument

# Locate the full sentence that needs partial re-formatting.
$full = $d.Content
$found = $full.Find.Execute(
    "包括但不限于以下：摘要 ；问题定义；技术现状；所采用或提出的方法；实验结果；结论；参考文献",
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if (-not $found) {
    throw "Target sentence not found"
}

$baseStart = $full.Start
$baseEnd = $full.End

# Offsets (character counts) within the matched sentence:
#   0  .. 9  -> "包括但不限于以下："           (stays unformatted)
#   9  .. 11 -> "摘要"                          (turns red + yellow highlight)
#   11 .. 12 -> " "                              (turns red + yellow highlight)
#   12 .. 41 -> "；问题定义；技术现状；所采用或提出的方法；实验结果；结论；" (turns red + yellow highlight)
#   41 .. 45 -> "参考文献"                      (stays unformatted)

$rHighlighted = $d.Range($baseStart + 9, $baseStart + 41)
$rHighlighted.Font.Color = 255
$rHighlighted.Font.HighlightColorIndex = 7
